$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.975105
$ws.Range("H2").Value = 56.925315
$ws.Range("I2").Value = 0.9552145540969871
$ws.Range("J2").Value = 0.955214554096987
$ws.Range("M2").Value = 0.2303363333333333
$ws.Range("N2").Value = 0.691009
$ws.Range("O2").Value = 0.0420565315194687
$ws.Range("P2").Value = 0.0420565315194687
$ws.Range("Q2").Value = 4.370656110315
$ws.Range("R2").Value = 39.335904992835
$ws.Range("S2").Value = 0.04017301100223517
$ws.Range("T2").Value = 0.04017301100223517
$ws.Range("G3").Value = 18.975105
$ws.Range("H3").Value = 56.925315
$ws.Range("I3").Value = 0.9552145540969871
$ws.Range("J3").Value = 0.955214554096987
$ws.Range("O3").Value = 0.8440851393264226
$ws.Range("P3").Value = 0.8440851393264227
$ws.Range("Q3").Value = 87.72016470534
$ws.Range("R3").Value = 789.48148234806
$ws.Range("S3").Value = 0.806282409981582
$ws.Range("T3").Value = 0.806282409981582
$ws.Range("G4").Value = 18.975105
$ws.Range("H4").Value = 56.925315
$ws.Range("I4").Value = 0.9552145540969871
$ws.Range("J4").Value = 0.955214554096987
$ws.Range("M4").Value = 0.6235823333333333
$ws.Range("N4").Value = 1.870747
$ws.Range("O4").Value = 0.1138583291541087
$ws.Range("P4").Value = 0.1138583291541087
$ws.Range("Q4").Value = 11.832540251145
$ws.Range("R4").Value = 106.492862260305
$ws.Range("S4").Value = 0.1087591331131699
$ws.Range("T4").Value = 0.1087591331131699
$ws.Range("I5").Value = 0.01570916103663723
$ws.Range("J5").Value = 0.01570916103663723
$ws.Range("M5").Value = 0.2303363333333333
$ws.Range("N5").Value = 0.691009
$ws.Range("O5").Value = 0.0420565315194687
$ws.Range("P5").Value = 0.0420565315194687
$ws.Range("Q5").Value = 0.07187844906488888
$ws.Range("R5").Value = 0.646906041584
$ws.Range("S5").Value = 0.0006606728262817433
$ws.Range("T5").Value = 0.0006606728262817433
$ws.Range("I6").Value = 0.01570916103663723
$ws.Range("J6").Value = 0.01570916103663723
$ws.Range("O6").Value = 0.8440851393264226
$ws.Range("P6").Value = 0.8440851393264227
$ws.Range("S6").Value = 0.01325986938231115
$ws.Range("T6").Value = 0.01325986938231115
$ws.Range("I7").Value = 0.01570916103663723
$ws.Range("J7").Value = 0.01570916103663723
$ws.Range("M7").Value = 0.6235823333333333
$ws.Range("N7").Value = 1.870747
$ws.Range("O7").Value = 0.1138583291541087
$ws.Range("P7").Value = 0.1138583291541087
$ws.Range("Q7").Value = 0.1945942714968889
$ws.Range("R7").Value = 1.751348443472
$ws.Range("S7").Value = 0.001788618828044341
$ws.Range("T7").Value = 0.001788618828044341
$ws.Range("G8").Value = 0.5775933333333333
$ws.Range("H8").Value = 1.73278
$ws.Range("I8").Value = 0.02907628486637583
$ws.Range("J8").Value = 0.02907628486637583
$ws.Range("M8").Value = 0.2303363333333333
$ws.Range("N8").Value = 0.691009
$ws.Range("O8").Value = 0.0420565315194687
$ws.Range("P8").Value = 0.0420565315194687
$ws.Range("Q8").Value = 0.1330407305577778
$ws.Range("R8").Value = 1.19736657502
$ws.Range("S8").Value = 0.001222847690951786
$ws.Range("T8").Value = 0.001222847690951786
$ws.Range("G9").Value = 0.5775933333333333
$ws.Range("H9").Value = 1.73278
$ws.Range("I9").Value = 0.02907628486637583
$ws.Range("J9").Value = 0.02907628486637583
$ws.Range("O9").Value = 0.8440851393264226
$ws.Range("P9").Value = 0.8440851393264227
$ws.Range("Q9").Value = 2.670160841413333
$ws.Range("R9").Value = 24.03144757272
$ws.Range("S9").Value = 0.0245428599625296
$ws.Range("T9").Value = 0.0245428599625296
$ws.Range("G10").Value = 0.5775933333333333
$ws.Range("H10").Value = 1.73278
$ws.Range("I10").Value = 0.02907628486637583
$ws.Range("J10").Value = 0.02907628486637583
$ws.Range("M10").Value = 0.6235823333333333
$ws.Range("N10").Value = 1.870747
$ws.Range("O10").Value = 0.1138583291541087
$ws.Range("P10").Value = 0.1138583291541087
$ws.Range("Q10").Value = 0.3601769985177777
$ws.Range("R10").Value = 3.24159298666
$ws.Range("S10").Value = 0.003310577212894449
$ws.Range("T10").Value = 0.003310577212894449
